$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the author name for the GWAS study row from "Converge" to "CONVERGE"
$ws.Range("A3").Value = "CONVERGE"

# Update the active selection to A4 (as reflected in the saved view state)
$ws.Range("A4").Select()
